$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-20 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-21 Friday", 2)

$d.Content.Find.Execute("942×3=2826", $true, $false, $false, $false, $false, $true, 1, $false, "374×9=3366", 2)
$d.Content.Find.Execute("911×5=4555", $true, $false, $false, $false, $false, $true, 1, $false, "949×8=7592", 2)
$d.Content.Find.Execute("176×2=352", $true, $false, $false, $false, $false, $true, 1, $false, "533×5=2665", 2)
$d.Content.Find.Execute("784×5=3920", $true, $false, $false, $false, $false, $true, 1, $false, "749×4=2996", 2)
$d.Content.Find.Execute("791×2=1582", $true, $false, $false, $false, $false, $true, 1, $false, "740×7=5180", 2)
$d.Content.Find.Execute("850×5=4250", $true, $false, $false, $false, $false, $true, 1, $false, "569×5=2845", 2)
$d.Content.Find.Execute("518×5=2590", $true, $false, $false, $false, $false, $true, 1, $false, "940×6=5640", 2)
$d.Content.Find.Execute("112×5=560", $true, $false, $false, $false, $false, $true, 1, $false, "475×4=1900", 2)
$d.Content.Find.Execute("633×7=4431", $true, $false, $false, $false, $false, $true, 1, $false, "108×2=216", 2)
$d.Content.Find.Execute("630×2=1260", $true, $false, $false, $false, $false, $true, 1, $false, "782×9=7038", 2)
$d.Content.Find.Execute("406×3=1218", $true, $false, $false, $false, $false, $true, 1, $false, "488×8=3904", 2)
$d.Content.Find.Execute("184×9=1656", $true, $false, $false, $false, $false, $true, 1, $false, "925×8=7400", 2)
$d.Content.Find.Execute("815×3=2445", $true, $false, $false, $false, $false, $true, 1, $false, "550×7=3850", 2)
$d.Content.Find.Execute("306×7=2142", $true, $false, $false, $false, $false, $true, 1, $false, "886×7=6202", 2)
$d.Content.Find.Execute("418×7=2926", $true, $false, $false, $false, $false, $true, 1, $false, "829×3=2487", 2)
$d.Content.Find.Execute("527×2=1054", $true, $false, $false, $false, $false, $true, 1, $false, "600×2=1200", 2)
$d.Content.Find.Execute("377×4=1508", $true, $false, $false, $false, $false, $true, 1, $false, "450×2=900", 2)
$d.Content.Find.Execute("268×5=1340", $true, $false, $false, $false, $false, $true, 1, $false, "471×5=2355", 2)
$d.Content.Find.Execute("308×6=1848", $true, $false, $false, $false, $false, $true, 1, $false, "202×3=606", 2)
$d.Content.Find.Execute("186×2=372", $true, $false, $false, $false, $false, $true, 1, $false, "620×3=1860", 2)
$d.Content.Find.Execute("337×2=674", $true, $false, $false, $false, $false, $true, 1, $false, "980×7=6860", 2)
$d.Content.Find.Execute("598×9=5382", $true, $false, $false, $false, $false, $true, 1, $false, "266×8=2128", 2)
$d.Content.Find.Execute("429×7=3003", $true, $false, $false, $false, $false, $true, 1, $false, "934×9=8406", 2)
$d.Content.Find.Execute("407×2=814", $true, $false, $false, $false, $false, $true, 1, $false, "958×4=3832", 2)
$d.Content.Find.Execute("339×9=3051", $true, $false, $false, $false, $false, $true, 1, $false, "667×5=3335", 2)
